$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 490, shifting existing rows 490-576 down to 493-579
$ws.Rows("490:492").Insert()

# Populate the newly inserted rows with the new data (boilerplate columns A,B,C,E,F,G,H,I,J,K,Q,T
# copied from the surrounding rows; only D,L,M,N,O,P,R,S differ as specified by the edit)

# Row 490
$ws.Range("A490").Value = 9
$ws.Range("B490").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C490").Value = "Metropolitana"
$ws.Range("D490").Value = 44504
$ws.Range("E490").Value = 13
$ws.Range("F490").Value = "Fruta"
$ws.Range("G490").Value = 100101
$ws.Range("H490").Value = "Berries"
$ws.Range("I490").Value = 100112025
$ws.Range("J490").Value = "Frutilla"
$ws.Range("K490").Value = "Sin especificar"
$ws.Range("L490").Value = "Especial"
$ws.Range("M490").Value = 980
$ws.Range("N490").Value = 5000
$ws.Range("O490").Value = 5000
$ws.Range("P490").Value = 5000
$ws.Range("Q490").Value = "`$/bandeja 7 kilos"
$ws.Range("R490").Value = "Provincia de San Antonio"
$ws.Range("S490").Value = 714
$ws.Range("T490").Value = 7

# Row 491
$ws.Range("A491").Value = 9
$ws.Range("B491").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C491").Value = "Metropolitana"
$ws.Range("D491").Value = 44504
$ws.Range("E491").Value = 13
$ws.Range("F491").Value = "Fruta"
$ws.Range("G491").Value = 100101
$ws.Range("H491").Value = "Berries"
$ws.Range("I491").Value = 100112025
$ws.Range("J491").Value = "Frutilla"
$ws.Range("K491").Value = "Sin especificar"
$ws.Range("L491").Value = "Primera"
$ws.Range("M491").Value = 1200
$ws.Range("N491").Value = 4000
$ws.Range("O491").Value = 4000
$ws.Range("P491").Value = 4000
$ws.Range("Q491").Value = "`$/bandeja 7 kilos"
$ws.Range("R491").Value = "Provincia de San Antonio"
$ws.Range("S491").Value = 571
$ws.Range("T491").Value = 7

# Row 492
$ws.Range("A492").Value = 9
$ws.Range("B492").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C492").Value = "Metropolitana"
$ws.Range("D492").Value = 44504
$ws.Range("E492").Value = 13
$ws.Range("F492").Value = "Fruta"
$ws.Range("G492").Value = 100101
$ws.Range("H492").Value = "Berries"
$ws.Range("I492").Value = 100112025
$ws.Range("J492").Value = "Frutilla"
$ws.Range("K492").Value = "Sin especificar"
$ws.Range("L492").Value = "Segunda"
$ws.Range("M492").Value = 950
$ws.Range("N492").Value = 3000
$ws.Range("O492").Value = 3000
$ws.Range("P492").Value = 3000
$ws.Range("Q492").Value = "`$/bandeja 7 kilos"
$ws.Range("R492").Value = "Provincia de San Antonio"
$ws.Range("S492").Value = 429
$ws.Range("T492").Value = 7
